# Insert a new data row above the current row 263 (weekly Ciboulette price
# record for Mercado Mayorista Lo Valledor de Santiago), pushing the old
# rows 263:284 down to 264:285. The inserted row carries the new week's
# price observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 263:284 down by one, duplicating row 263's formatting onto the
# freshly inserted blank row (mirrors Excel's default "insert" behaviour).
$ws.Rows.Item(263).Insert()

# Populate the newly inserted row 263 with the new observation.
$ws.Range("A263").Value = 6
$ws.Range("B263").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C263").Value = "Metropolitana"
$ws.Range("D263").Value = 44461
$ws.Range("E263").Value = 13
$ws.Range("F263").Value = 100112039
$ws.Range("G263").Value = "Ciboulette"
$ws.Range("H263").Value = "Sin especificar"
$ws.Range("I263").Value = "Primera"
$ws.Range("J263").Value = 830
$ws.Range("K263").Value = 1000
$ws.Range("L263").Value = 1200
$ws.Range("M263").Value = 1084
$ws.Range("N263").Value = "$/docena de atados"
$ws.Range("O263").Value = "Región Metropolitana"
$ws.Range("P263").Value = 361
$ws.Range("Q263").Value = 3
$ws.Range("R263").Value = "Hortaliza"
